# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match newly scraped counts.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" --------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 6598
$ws1.Range("F6").Value  = 80
$ws1.Range("F9").Value  = 5981
$ws1.Range("F10").Value = 40
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 1238
$ws1.Range("F15").Value = 384
$ws1.Range("F18").Value = 346
$ws1.Range("F21").Value = 4299
$ws1.Range("F22").Value = 42
$ws1.Range("F24").Value = 186
$ws1.Range("F25").Value = 18

# ---- Sheet "全部类型" -----------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 6598
$ws4.Range("F6").Value  = 80
$ws4.Range("F9").Value  = 5981
$ws4.Range("F10").Value = 40
$ws4.Range("F12").Value = 1238
$ws4.Range("F15").Value = 384
$ws4.Range("F18").Value = 346
$ws4.Range("F21").Value = 4299
$ws4.Range("F23").Value = 42
$ws4.Range("F25").Value = 186
$ws4.Range("F26").Value = 18
